$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: measure -> dimension renames for "rama" and "rama-descripcion"
$ws.Range("A3").Value = "iaest-dimension:rama"
$ws.Range("G3").Value = "iaest-dimension:rama-descripcion"

# Row 4: reclassify rama / comarca-nombre / ano / rama-descripcion as "dim",
# keep vab as "medida"
$ws.Range("A4").Value = "dim"
$ws.Range("B4").Value = "medida"
$ws.Range("C4").Value = "dim"
$ws.Range("F4").Value = "dim"
$ws.Range("G4").Value = "dim"

# Row 5: "rama" and "rama-descripcion" now typed as skos:Concept instead of xsd:string
$ws.Range("A5").Value = "skos:Concept"
$ws.Range("G5").Value = "skos:Concept"

# Row 6 (new): mapping file references for rama / rama-descripcion
$ws.Range("A6").Value = "mapping-rama.xlsx"
$ws.Range("G6").Value = "mapping-rama-descripcion.xlsx"
